$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.147.60"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "'1.656.18"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'218.11"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'0.5296"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.06331"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'20.44"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'4.520"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "'1.656.36"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'1.883.52"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'0.5494"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'0.0₅8214"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'26.132.87"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'4.604"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'191.36"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "'6.019"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'144.94"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'7.214"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "'1.468"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").Value = "'0.05745"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'3.563"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "'3.273"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "'1.601"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").Value = "'2.803"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "'0.9510"
$ws.Range("D37").Value = "'2.418"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "'0.5734"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'0.01609"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.806"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8508"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'104.47"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.044.25"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.004"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'1.797.70"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'56.91"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "'0.4346"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05155"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.845"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'1.443"
$ws.Range("E51").Value = "  -2.35%  "
